# Apply the edits described in the diff:
# 1. Rename the "MODEL_CONDITION" header text to "MODELCONDITION".
# 2. Drop the original column A (the 1/3/13/14/19 index column with the bordered
#    style) and shift every other column one place to the left, so that the data
#    which used to live in columns B:F now lives in columns A:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text before the column shift so it's easy to target by its
# well known (pre-shift) address.
$ws.Range("E1").Value = "MODELCONDITION"

# Remove the old column A entirely; Excel shifts B:F left into A:E for us.
$ws.Columns.Item(1).Delete()

# Make sure the sheet's used range / dimension reflects the new, narrower
# layout (A1:E6).
$ws.UsedRange | Out-Null
